$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated crypto price/volume data
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "29.915.93"
$ws.Range("E2").Value = "  +1.70%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.937.13"
$ws.Range("E3").Value = "  +1.32%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.007"
$ws.Range("E4").Value = "  -0.35%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "336.00"
$ws.Range("E5").Value = "  +3.12%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.005"
$ws.Range("E6").Value = "  -0.52%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4834"
$ws.Range("E7").Value = "  +0.46%  "
$ws.Range("E8").Value = "  +1.19%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.08167"
$ws.Range("E9").Value = "  -0.44%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.015"
$ws.Range("E10").Value = "  -0.36%  "
$ws.Range("E11").Value = "  +0.93%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.941.83"
$ws.Range("E12").Value = "  +1.12%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "6.088"
$ws.Range("E13").Value = "  +0.93%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.292"
$ws.Range("E14").Value = "  +1.10%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "91.18"
$ws.Range("E15").Value = "  +0.26%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.06847"
$ws.Range("E16").Value = "  +0.59%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.007"
$ws.Range("E17").Value = "  -0.40%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.00001035"
$ws.Range("E18").Value = "  -0.30%  "
$ws.Range("E19").Value = "  +0.68%  "
$ws.Range("E20").Value = "  -0.47%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "29.904.53"
$ws.Range("E21").Value = "  +1.52%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.631"
$ws.Range("E22").Value = "  +0.25%  "
$ws.Range("E23").Value = "  +1.29%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.180"
$ws.Range("E24").Value = "  -0.86%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.196.77"
$ws.Range("E25").Value = "  +2.13%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "6.698"
$ws.Range("E26").Value = "  +1.68%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "156.81"
$ws.Range("E27").Value = "  -0.09%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "20.02"
$ws.Range("E28").Value = "  -0.24%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.097"
$ws.Range("E29").Value = "  -0.35%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "121.63"
$ws.Range("E30").Value = "  +1.34%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.007"
$ws.Range("E31").Value = "  -1.21%  "
$ws.Range("E32").Value = "  +0.87%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.581"
$ws.Range("E33").Value = "  +0.62%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.420"
$ws.Range("E34").Value = "  +4.00%  "
$ws.Range("E35").Value = "  -0.42%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.06559"
$ws.Range("E36").Value = "  +7.41%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02286"
$ws.Range("E37").Value = "  +0.17%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.207"
$ws.Range("E38").Value = "  +2.40%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.5968"
$ws.Range("E39").Value = "  -0.04%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "7.980"
$ws.Range("E40").Value = "  -0.71%  "
$ws.Range("E41").Value = "  -0.44%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1847"
$ws.Range("E42").Value = "  -0.16%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.491"
$ws.Range("E43").Value = "  +3.01%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.269"
$ws.Range("E44").Value = "  +1.89%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "12.39"
$ws.Range("E45").Value = "  +0.06%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.07483"
$ws.Range("E46").Value = "  -1.35%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5570"
$ws.Range("E47").Value = "  +0.02%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.990"
$ws.Range("E48").Value = "  +2.09%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "117.06"
$ws.Range("E49").Value = "  -0.17%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "72.63"
$ws.Range("E50").Value = "  +0.25%  "
$ws.Range("E51").Value = "  -0.35%  "
